$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo/extra word in the exception's system-response text
# (cell D14, row for "Exceção 1 [Dados inválidos] (Passo 5)")
$ws.Range("D14").Value = "5.1. Informa que as credenciais são inválidas"

# Reset the view: scroll back to the top-left and move the active
# selection to D15 (matches the author's cursor position when they
# saved after writing a bit more of the report)
$ws.Range("A1").Select() | Out-Null
$ws.Range("D15").Select() | Out-Null
